# Update countries & provincias Spain
# - Reorder several country names in the "Pais" sheet (rows 206-214) to
#   reflect the latest rankings, carrying their active-case (D) / deaths
#   (H) values along with the correct country.
# - Bump the "last updated" timestamp in A1 from 11:27 to 11:28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Countries around rows 206-214 swap order -----------------------------
# Row 206: Islas Malvinas -> Groenlandia (values unchanged: 13 / 13 / 0)
$ws.Cells.Item(206, 1).Value = "Groenlandia"

# Row 207: Groenlandia -> Islas Malvinas (values unchanged: 13 / 13 / 0)
$ws.Cells.Item(207, 1).Value = "Islas Malvinas"

# Row 208: Islas Turcas y Caicos -> Santa Sede (D 11->12, H 1->0)
$ws.Cells.Item(208, 1).Value = "Santa Sede"
$ws.Cells.Item(208, 4).Value = 12
$ws.Cells.Item(208, 8).Value = 0

# Row 209: Santa Sede -> Islas Turcas y Caicos (D 12->11, H 0->1)
$ws.Cells.Item(209, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(209, 4).Value = 11
$ws.Cells.Item(209, 8).Value = 1

# Row 210: Seychelles -> Montserrat (D 11->10, H 0->1)
$ws.Cells.Item(210, 1).Value = "Montserrat"
$ws.Cells.Item(210, 4).Value = 10
$ws.Cells.Item(210, 8).Value = 1

# Row 211: Montserrat -> Seychelles (D 10->11, H 1->0)
$ws.Cells.Item(211, 1).Value = "Seychelles"
$ws.Cells.Item(211, 4).Value = 11
$ws.Cells.Item(211, 8).Value = 0

# Row 212: Sahara Occidental stays the same (no change)

# Row 213: Papua Nueva Guinea -> Islas Virgenes Britanicas (D 8->7, H 0->1)
$ws.Cells.Item(213, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(213, 4).Value = 7
$ws.Cells.Item(213, 8).Value = 1

# Row 214: Islas Virgenes Britanicas -> Papua Nueva Guinea (D 7->8, H 1->0)
$ws.Cells.Item(214, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(214, 4).Value = 8
$ws.Cells.Item(214, 8).Value = 0

# Row 215: Bonaire, San Eustaquio y Saba stays the same (no change)

# --- Bump the "updated at" timestamp --------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 13 de Junio de 2020 a las 11:28"
